# tradexcb_strategy.xlsx update
# - upgrade to version v0.1.8
# - update transaction types to variables at required positions
#   (expiry/instrument refreshed for the existing Buy leg, a new Sell leg
#   added for the same instrument/expiry)
# - corrected mapping logic for order status (vwap / moving_average flags)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Buy leg) ---------------------------------
# refresh expiry/instrument to the new series
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2022-06-09"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "NIFTY2260916500CE"

# corrected mapping logic for order status
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "1"
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = "YES"
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = "NO"

# --- Add new row 3 (Sell leg) ----------------------------------------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "Sell"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.0"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0.0"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "MARKET"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "MIS"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "5.0"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "NFO"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "NIFTY"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2022-06-09"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "NIFTY2260916500CE"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "0"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "NO"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "0"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "NO"
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "Value"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "7"
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "Value"
$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "2"
$ws.Range("S3").NumberFormat = "@"
$ws.Range("S3").Value = "Value"
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = "15"
$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value = "1"
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value = "YES"
$ws.Range("W3").NumberFormat = "@"
$ws.Range("W3").Value = "existing"
$ws.Range("X3").NumberFormat = "@"
$ws.Range("X3").Value = "NO"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "21"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "3"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "new"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "21"
$ws.Range("AC3").NumberFormat = "@"
$ws.Range("AC3").Value = "NO"
$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AD3").Value = "new"
$ws.Range("AE3").NumberFormat = "@"
$ws.Range("AE3").Value = "Default"
